$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Elemento")
Write-Host "Sheet name: " $ws.Name
Write-Host "A1: " $ws.Range("A1").Value
